# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-24, columns C:G) is re-sorted so the
# records are grouped by period (1909..1912) instead of by worker, and the
# "Valor Mora" (column G) figures for YULY PAULINA MUÑOZ OSPINO are updated
# from 743000 to 828116.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Worker master data (Tipo Doc column B is "CC" everywhere and is untouched).
$jhon  = @{ Doc = "73213301";   Nombre = "JHON JAIRO BARRIOS SEPULVEDA" }
$belkis = @{ Doc = "1047403776"; Nombre = "BELKIS MARIA FLOREZ GONZALEZ" }
$yuly  = @{ Doc = "1128048125"; Nombre = "YULY PAULINA MUÑOZ OSPINO" }

# New row order/content for rows 16-24: (worker, periodo, valorMora, salario)
$rows = @(
    @{ Row = 16; Worker = $jhon;   Periodo = "1909"; Mora = 40000; Salario = 1000000 },
    @{ Row = 17; Worker = $yuly;   Periodo = "1909"; Mora = 33125; Salario = 828116 },
    @{ Row = 18; Worker = $jhon;   Periodo = "1910"; Mora = 40000; Salario = 1000000 },
    @{ Row = 19; Worker = $yuly;   Periodo = "1910"; Mora = 33125; Salario = 828116 },
    @{ Row = 20; Worker = $jhon;   Periodo = "1911"; Mora = 40000; Salario = 1000000 },
    @{ Row = 21; Worker = $belkis; Periodo = "1911"; Mora = 18771; Salario = 878000 },
    @{ Row = 22; Worker = $yuly;   Periodo = "1911"; Mora = 33125; Salario = 828116 },
    @{ Row = 23; Worker = $jhon;   Periodo = "1912"; Mora = 40000; Salario = 1000000 },
    @{ Row = 24; Worker = $belkis; Periodo = "1912"; Mora = 33125; Salario = 878000 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("C$n").Value = $r.Worker.Doc
    $ws.Range("D$n").Value = $r.Worker.Nombre
    $ws.Range("E$n").Value = $r.Periodo
    $ws.Range("F$n").Value = $r.Mora
    $ws.Range("G$n").Value = $r.Salario
}
